$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.11"
$ws.Range("E2").Value = "'-2.22%"
$ws.Range("D3").Value = "'35.89"
$ws.Range("E3").Value = "'1.15%"
$ws.Range("D4").Value = "'5.053"
$ws.Range("E4").Value = "'-1.33%"
$ws.Range("D5").Value = "'0.08078"
$ws.Range("E5").Value = "'-1.53%"
$ws.Range("D6").Value = "'1.946"
$ws.Range("E6").Value = "'-4.88%"
$ws.Range("D7").Value = "'7.811"
$ws.Range("E7").Value = "'-1.91%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9287"
$ws.Range("E8").Value = "'0.08%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1495"
$ws.Range("E9").Value = "'37.84%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1907"
$ws.Range("E10").Value = "'-1.39%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08996"
$ws.Range("E11").Value = "'-6.36%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03478"
$ws.Range("E12").Value = "'-3.76%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09849"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001412"
$ws.Range("E14").Value = "'-1.24%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005877"
$ws.Range("E15").Value = "'2.25%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.529"
$ws.Range("E16").Value = "'1.56%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.069"
$ws.Range("E17").Value = "'-1.35%"
$ws.Range("E18").Value = "'-0.70%"
$ws.Range("D19").Value = "'0.3452"
$ws.Range("E19").Value = "'1.01%"
$ws.Range("D20").Value = "'0.1286"
$ws.Range("E20").Value = "'-1.02%"
$ws.Range("D21").Value = "'5.028"
$ws.Range("E21").Value = "'-1.36%"
$ws.Range("E22").Value = "'9.38%"
$ws.Range("D23").Value = "'0.04482"
$ws.Range("E23").Value = "'-1.34%"
$ws.Range("D24").Value = "'0.001203"
$ws.Range("E24").Value = "'-1.83%"
$ws.Range("D25").Value = "'0.004815"
$ws.Range("D26").Value = "'0.0001226"
$ws.Range("E26").Value = "'-2.01%"
$ws.Range("E27").Value = "'-32.21%"
$ws.Range("D39").Value = "'0.01902"
$ws.Range("E39").Value = "'-3.78%"
$ws.Range("D40").Value = "'0.04784"
$ws.Range("E40").Value = "'-2.70%"
$ws.Range("D41").Value = "'0.01057"
$ws.Range("E41").Value = "'10.08%"
$ws.Range("D42").Value = "'0.007324"
$ws.Range("E42").Value = "'-6.66%"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'-2.65%"
$ws.Range("D44").Value = "'0.002104"
$ws.Range("E44").Value = "'-0.64%"
$ws.Range("D45").Value = "'0.01077"
$ws.Range("E45").Value = "'-6.87%"
$ws.Range("D46").Value = "'0.00006111"
$ws.Range("E46").Value = "'-5.66%"
$ws.Range("D47").Value = "'0.00000000747"
$ws.Range("E47").Value = "'-0.56%"
$ws.Range("E48").Value = "'-3.00%"
$ws.Range("E49").Value = "'27.66%"
$ws.Range("D50").Value = "'0.00002091"
$ws.Range("E50").Value = "'-0.56%"
$ws.Range("D51").Value = "'0.0001991"
$ws.Range("E51").Value = "'-0.56%"
